$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (E2:T2)
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.319712
$ws.Range("H2").Value = 3.959136
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.08121833333333332
$ws.Range("N2").Value = 0.243655
$ws.Range("O2").Value = 0.9961121149930909
$ws.Range("P2").Value = 0.9961121149930909
$ws.Range("Q2").Value = 0.10718480912
$ws.Range("R2").Value = 0.96466328208
$ws.Range("S2").Value = 0.9961121149930909
$ws.Range("T2").Value = 0.9961121149930909

# Row 3 (E3:T3)
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.319712
$ws.Range("H3").Value = 3.959136
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.000317
$ws.Range("N3").Value = 0.000951
$ws.Range("O3").Value = 0.003887885006909071
$ws.Range("P3").Value = 0.00388788500690907
$ws.Range("Q3").Value = 0.000418348704
$ws.Range("R3").Value = 0.003765138336
$ws.Range("S3").Value = 0.003887885006909071
$ws.Range("T3").Value = 0.00388788500690907
